$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (8th column) to hold the CO2/(CO+CO2) ratio
$ws.Columns.Item(8).Insert()

# Header
$ws.Cells.Item(1, 8).Value = "CO2/(CO+CO2)"

# Formula for each data row: CO2/(CO+CO2) = F/(E+F)
for ($r = 2; $r -le 29; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Formula = "=F$r/(E$r+F$r)"
    $cell.NumberFormat = "General"
    $cell.HorizontalAlignment = $xlGeneral
    $cell.VerticalAlignment = $xlBottom
    $cell.WrapText = $false
}

$ws.Range("H2").Select()
